$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "43.144.07"
$ws.Range("E2").Value = "  +1.68%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.382.38"
$ws.Range("E3").Value = "  +3.95%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.59%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.31%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.40%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.16%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.501"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.93%  "

# Row 10 - Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.04%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +1.15%  "

# Row 12 - was Chainlink, now TRON
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.122"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.46%  "

# Row 13 - was TRON, now Chainlink
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.51%  "

# Row 14 - Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.79"
$ws.Range("D14").Style = "Normal"

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.749.96"
$ws.Range("E15").Value = "  +3.81%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.370.16"
$ws.Range("E16").Value = "  +2.83%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  +3.94%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "43.099.96"
$ws.Range("E18").Value = "  +1.68%  "

# Row 19 - InternetComputer(DFINITY)
$ws.Range("E19").Value = "  +0.36%  "

# Row 20 - Uniswap
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.13%  "

# Row 21 - ShibaInu
$ws.Range("D21").Value = "0.0₃0890"
$ws.Range("E21").Value = "  +0.28%  "

# Row 22 - Litecoin
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.10%  "

# Row 23 - BitcoinCash
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.46%  "

# Row 24 - ImmutableX
$ws.Range("E24").Value = "  -2.33%  "

# Row 25 - was PancakeSwap, now Dai
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.13%  "

# Row 26 - was Dai, now PancakeSwap
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.07%  "

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.16%  "

# Row 28 - Toncoin
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.72%  "

# Row 29 - Cosmos
$ws.Range("E29").Value = "  +1.06%  "

# Row 30 - InjectiveProtocol
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.85%  "

# Row 31 - FirstDigitalUSD
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.12%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +2.39%  "

# Row 33 - Hedera
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0733"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.25%  "

# Row 34 - Celestia
$ws.Range("E34").Value = "  -1.48%  "

# Row 35 - ARBITRUM
$ws.Range("E35").Value = "  +7.03%  "

# Row 36 - RenderToken
$ws.Range("E36").Value = "  -0.74%  "

# Row 37 - Kaspa
$ws.Range("E37").Value = "  +1.97%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  +4.21%  "

# Row 40 - EnergySwap
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.78%  "

# Row 41 - Stellar
$ws.Range("E41").Value = "  +0.33%  "

# Row 42 - Monero
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -33.93%  "

# Row 43 - Maker
$ws.Range("D43").Value = "1.956.95"
$ws.Range("E43").Value = "  -0.19%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  +0.74%  "

# Row 45 - ApeXProtocol
$ws.Range("E45").Value = "  +2.23%  "

# Row 46 - NEARProtocol
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.38%  "

# Row 47 - FraxShare
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.92%  "

# Row 48 - RocketPoolETH
$ws.Range("D48").Value = "2.605.61"
$ws.Range("E48").Value = "  +3.47%  "

# Row 49 - MultiversX
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.50%  "

# Row 50 - Stacks
$ws.Range("E50").Value = "  +2.04%  "

# Row 51 - BitcoinSV
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.88%  "
